$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132781267166138
$ws.Range("B1").Value = 1.298191905021667
$ws.Range("C1").Value = 1.044164061546326
$ws.Range("D1").Value = 5.361801624298096
$ws.Range("E1").Value = 1.814984917640686
